# Generate Report for Handback
#
# The 66a3607c-7286-46b5-b04e-156b1d1d137f.md file has now been handed back
# in sync with en-US, so update its status (was "Ready for handoff") and
# record the new handback timestamps, clearing the stale version-mismatch
# error message.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 66a3607c-...md is row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row for 66a3607c-...md is row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-09-07 09:03:17"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: row for 66a3607c-...md is row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-09-07 09:03:38"
$wsDeDe.Range("P3").Value = ""
